$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.856.05"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "1.635.14"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'0.5019"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'0.2570"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").Value = "'0.06404"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'19.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "'0.07701"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.242"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.634.16"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "1.861.11"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'0.5435"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "0.0₅7922"
$ws.Range("D17").Value = "'63.45"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "25.881.95"
$ws.Range("E18").Value = "  -1.19%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "'203.07"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "'4.325"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").Value = "'9.938"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "'5.976"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'1.920"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +11.04%  "
$ws.Range("D26").Value = "'140.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "'0.1143"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'6.703"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("D31").Value = "'0.04994"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'3.258"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "'3.176"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "'1.539"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "'2.366"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "1.168.88"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("D37").Value = "'0.8935"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("D38").Value = "'2.615"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("D39").Value = "'0.5612"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "'0.01558"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").Value = "'2.556"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "'5.672"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'0.8080"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "'99.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "1.772.78"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("D48").Value = "'0.4516"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'54.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "'0.05085"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.63%  "
